# Natmi following Dr Hou advice
# Rewrite the data rows (2-10) of Sheet1 with updated Ntf3-Ntrk2 signalling
# results that now include the "sCs" cluster as both a sending and a
# target cluster (3x3 grid instead of 2x3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Ntf3","Ntrk2","ECs",3,1,8.859944,26.579832,0.4921128329655918,0.4921128329655918,3,1,0.5314143333333333,1.594243,0.01762280440365445,0.01762280440365445,4.708301234130666,42.374711107176,0.008672408199880898,0.008672408199880896),
    @("ECs","Ntf3","Ntrk2","FAPs",3,1,8.859944,26.579832,0.4921128329655918,0.4921128329655918,3,1,24.72809866666667,74.184296,0.8200351754599551,0.820035175459955,219.0895694131414,1971.806124718272,0.4035498333270346,0.4035498333270345),
    @("ECs","Ntf3","Ntrk2","sCs",3,1,8.859944,26.579832,0.4921128329655918,0.4921128329655918,3,1,4.895411333333333,14.686234,0.1623420201363905,0.1623420201363905,43.37307027029866,390.357632432688,0.07989059143867629,0.07989059143867629),
    @("FAPs","Ntf3","Ntrk2","ECs",3,1,5.867977666666667,17.603933,0.3259283708025871,0.3259283708025871,3,1,0.5314143333333333,1.594243,0.01762280440365445,0.01762280440365445,3.118327439746556,28.064946957719,0.005743771928255752,0.005743771928255751),
    @("FAPs","Ntf3","Ntrk2","FAPs",3,1,5.867977666666667,17.603933,0.3259283708025871,0.3259283708025871,3,1,24.72809866666667,74.184296,0.8200351754599551,0.820035175459955,145.1039307151298,1305.935376436168,0.2672727287384768,0.2672727287384767),
    @("FAPs","Ntf3","Ntrk2","sCs",3,1,5.867977666666667,17.603933,0.3259283708025871,0.3259283708025871,3,1,4.895411333333333,14.686234,0.1623420201363905,0.1623420201363905,28.72616437314689,258.535479358322,0.05291187013585454,0.05291187013585454),
    @("sCs","Ntf3","Ntrk2","ECs",2,0.6666666666666666,3.275965666666667,9.827897,0.1819587962318212,0.1819587962318212,3,1,0.5314143333333333,1.594243,0.01762280440365445,0.01762280440365445,1.740895110774555,15.668055996971,0.003206624275517801,0.0032066242755178),
    @("sCs","Ntf3","Ntrk2","FAPs",2,0.6666666666666666,3.275965666666667,9.827897,0.1819587962318212,0.1819587962318212,3,1,24.72809866666667,74.184296,0.8200351754599551,0.820035175459955,81.00840223394577,729.075620105512,0.1492126133944437,0.1492126133944437),
    @("sCs","Ntf3","Ntrk2","sCs",2,0.6666666666666666,3.275965666666667,9.827897,0.1819587962318212,0.1819587962318212,3,1,4.895411333333333,14.686234,0.1623420201363905,0.1623420201363905,16.03719945221088,144.334795069898,0.02953955856185969,0.02953955856185969)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}

$wb.Save()
